# Excel DSA search features: rename the INDEX column header to ID.
# (see commit message: "Added methods to search for attribute values
# (ID, FIRST_NAME, LAST_NAME, EMAIL_ADDRESS, BUSINESS_NAME)")
#
# The workbook's only worksheet (DATA_USER_FORM) has a header row whose
# first column was labelled "INDEX"; it becomes "ID" so it lines up with
# the new ID-based lookup helpers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ID"
